# Apply updated coin data from the 2022-12-21 14:24:54 UTC GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'248.90"
$ws.Range("D2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'22.75"
$ws.Range("D3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = "'5.271"
$ws.Range("D4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'0.05699"
$ws.Range("D5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'3.409"
$ws.Range("D6").Style = "Normal"
# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = "'0.8054"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '6MXTokenMX'
# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = "'0.9038"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '7FTXTokenFTT'
# Row 9
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = "'0.01111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8OneONE'
# Row 10
$ws.Range("D10").Value = "'0.1422"
$ws.Range("D10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.07445"
$ws.Range("D11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'0.03083"
$ws.Range("D12").Style = "Normal"
# Row 13
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D13").Value = "'0.1292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12ProBitTokenPROB'
# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.03002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitrueCoinBTR'
# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09385"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14BitMartTokenBMX'
# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = "'3.858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15MCDexMCB'
# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = "'0.001581"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16BitForexTokenBF'
# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = "'0.04804"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17CoinExTokenCET'
# Row 19
$ws.Range("B19").Value = 'UpBots'
$ws.Range("C19").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D19").Value = "'0.01827"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '18UpBotsUBXTBestin24h'
# Row 20
$ws.Range("D20").Value = "'0.006423"
$ws.Range("D20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = "'0.004992"
$ws.Range("D21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'0.0009982"
$ws.Range("D22").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'3.696"
$ws.Range("D24").Style = "Normal"
# Row 25
$ws.Range("B25").Value = 'KuCoinToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D25").Value = "'6.355"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '24KuCoinTokenKCS'
# Row 26
$ws.Range("B26").Value = 'BTSEToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D26").Value = "'2.201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '25BTSETokenBTSE'
# Row 27
$ws.Range("B27").Value = 'BitpandaEcosystemToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D27").Value = "'0.3300"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26BitpandaEcosystemTokenBEST'
# Row 40
$ws.Range("D40").Value = "'0.03983"
$ws.Range("D40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.006810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICK'
# Row 42
$ws.Range("D42").Value = "'0.1069"
$ws.Range("D42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = "'0.002771"
$ws.Range("D43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'0.007707"
$ws.Range("D44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'0.00005593"
$ws.Range("D45").Style = "Normal"
# Row 47
$ws.Range("D47").Value = "'0.4990"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
# Row 48
$ws.Range("D48").Value = "'0.2017"
$ws.Range("D48").Style = "Normal"
